$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header for column B ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# Apply column-A's existing date style (from A2) down through A22 before
# filling in the new rows, so every date cell keeps the same formatting.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column A: annual (year-end) date series, row 2 through row 22
$ws.Cells.Item(2, 1).Value = 38717
$ws.Cells.Item(3, 1).Value = 39082
$ws.Cells.Item(4, 1).Value = 39447
$ws.Cells.Item(5, 1).Value = 39813
$ws.Cells.Item(6, 1).Value = 40178
$ws.Cells.Item(7, 1).Value = 40543
$ws.Cells.Item(8, 1).Value = 40908
$ws.Cells.Item(9, 1).Value = 41274
$ws.Cells.Item(10, 1).Value = 41639
$ws.Cells.Item(11, 1).Value = 42004
$ws.Cells.Item(12, 1).Value = 42369
$ws.Cells.Item(13, 1).Value = 42735
$ws.Cells.Item(14, 1).Value = 43100
$ws.Cells.Item(15, 1).Value = 43465
$ws.Cells.Item(16, 1).Value = 43830
$ws.Cells.Item(17, 1).Value = 44196
$ws.Cells.Item(18, 1).Value = 44561
$ws.Cells.Item(19, 1).Value = 44926
$ws.Cells.Item(20, 1).Value = 45291
$ws.Cells.Item(21, 1).Value = 45657
$ws.Cells.Item(22, 1).Value = 46022

# Column B: values shift down one row (old B2 becomes B3), with new rows
# appended through B21; B2 and B22 are left blank.
$ws.Cells.Item(2, 2).ClearContents()
$ws.Cells.Item(3, 2).Value = 6.681483765882756
$ws.Cells.Item(4, 2).Value = 5.732148352530309
$ws.Cells.Item(5, 2).Value = 6.181322443148352
$ws.Cells.Item(6, 2).Value = 9.97031398925483
$ws.Cells.Item(7, 2).Value = 5.968279190641868
$ws.Cells.Item(8, 2).Value = 6.387913216057295
$ws.Cells.Item(9, 2).Value = 3.489647115587391
$ws.Cells.Item(10, 2).Value = 3.062667370145955
$ws.Cells.Item(11, 2).Value = 1.40861416720266
$ws.Cells.Item(12, 2).Value = 1.975538030067248
$ws.Cells.Item(13, 2).Value = 2.149250550875026
$ws.Cells.Item(14, 2).Value = 2.516312190944614
$ws.Cells.Item(15, 2).Value = 2.494967260739056
$ws.Cells.Item(16, 2).Value = 0.9136132777513017
$ws.Cells.Item(17, 2).Value = 2.275661779503824
$ws.Cells.Item(18, 2).Value = 3.546865287857126
$ws.Cells.Item(19, 2).Value = 4.646251873334628
$ws.Cells.Item(20, 2).Value = 2.244754177395403
$ws.Cells.Item(21, 2).Value = 1.608247521160311
$ws.Cells.Item(22, 2).ClearContents()
